$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to remain Text so numeric-looking strings
# (e.g. "0.7498", "1.001") are not silently reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.335.87"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.935.83"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "0.7498"
$ws.Range("E5").Value = "  +5.27%  "
$ws.Range("D6").Value = "243.27"
$ws.Range("E6").Value = "  -2.12%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "28.04"
$ws.Range("E8").Value = "  +2.26%  "
$ws.Range("D9").Value = "0.3183"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").Value = "0.07033"
$ws.Range("E10").Value = "  -0.41%  "
$ws.Range("D11").Value = "0.7820"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").Value = "0.08049"
$ws.Range("D13").Value = "1.928.27"
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("D14").Value = "5.403"
$ws.Range("E14").Value = "  +0.45%  "
$ws.Range("D15").Value = "93.06"
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("D16").Value = "14.51"
$ws.Range("E16").Value = "  -0.62%  "
$ws.Range("D17").Value = "30.338.91"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("D18").Value = "6.133"
$ws.Range("E18").Value = "  +6.58%  "
$ws.Range("D19").Value = "252.39"
$ws.Range("E19").Value = "  -2.31%  "
$ws.Range("D20").Value = "0.000007988"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("D21").Value = "2.194.07"
$ws.Range("E21").Value = "  +0.44%  "
$ws.Range("D22").Value = "0.9998"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "6.706"
$ws.Range("E24").Value = "  -2.04%  "
$ws.Range("D25").Value = "9.573"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").Value = "165.28"
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").Value = "19.10"
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("D28").Value = "0.1303"
$ws.Range("E28").Value = "  +3.73%  "
$ws.Range("D29").Value = "2.202"
$ws.Range("E29").Value = "  -2.67%  "
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").Value = "1.546"
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("D32").Value = "4.439"
$ws.Range("E32").Value = "  +1.11%  "
$ws.Range("D33").Value = "4.148"
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("D34").Value = "0.05293"
$ws.Range("E34").Value = "  +3.09%  "
$ws.Range("E35").Value = "  +5.41%  "
$ws.Range("D36").Value = "0.7572"
$ws.Range("D37").Value = "2.784"
$ws.Range("E37").Value = "  +0.70%  "
$ws.Range("D38").Value = "0.01962"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").Value = "2.802"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").Value = "79.01"
$ws.Range("E40").Value = "  +2.05%  "
$ws.Range("D41").Value = "6.521"
$ws.Range("E41").Value = "  +2.57%  "
$ws.Range("D42").Value = "0.4518"
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("D43").Value = "1.983"
$ws.Range("E43").Value = "  -0.19%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").Value = "0.8403"
$ws.Range("D46").Value = "7.724"
$ws.Range("E46").Value = "  +4.08%  "
$ws.Range("D47").Value = "10.00"
$ws.Range("E47").Value = "  +2.90%  "
$ws.Range("D48").Value = "101.70"
$ws.Range("E48").Value = "  +1.03%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "37.66"
$ws.Range("E49").Value = "  +3.04%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.1249"
$ws.Range("E50").Value = "  +9.96%  "
$ws.Range("D51").Value = "966.79"
$ws.Range("E51").Value = "  +5.58%  "

# Restore default (unstyled) cell style now that values are written,
# so no stray style index is left attached to the edited cells.
$ws.Range("D2:E51").Style = "Normal"

